$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook's table (rows 16-19) lists workers with duplicate rows for
# "1050973054 / NATALY GUERRERO JIMENEZ" (row 17) and
# "1143366762 / YUNAY YICETH CASTILLO MOJICA" (row 18).
# The author's edit swaps these two workers' ID/Name between the two rows
# (the underlying shared-strings table was reordered while the cells that
# reference it stayed put, which has the net visual effect of exchanging
# the two rows' C/D values).

$c17 = $ws.Range("C17").Value2
$d17 = $ws.Range("D17").Value2
$c18 = $ws.Range("C18").Value2
$d18 = $ws.Range("D18").Value2

$ws.Range("C17").Value = $c18
$ws.Range("D17").Value = $d18
$ws.Range("C18").Value = $c17
$ws.Range("D18").Value = $d17
